$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = "reabertura shopping"
$ws.Cells.Item(2,3).Value = 0.6843937662750245
$ws.Cells.Item(3,2).Value = "reabertura academias"
$ws.Cells.Item(3,3).Value = 0.6825796157480096
$ws.Cells.Item(4,3).Value = 0.6866396568106133
$ws.Cells.Item(5,2).Value = "oximetro"
$ws.Cells.Item(5,3).Value = 0.5744918018441746
$ws.Cells.Item(6,2).Value = "oximetro de dedo"
$ws.Cells.Item(6,3).Value = 0.5631017583420228
$ws.Cells.Item(7,2).Value = "ivermectina coronavírus"
$ws.Cells.Item(7,3).Value = 0.5536793810526284
$ws.Cells.Item(8,2).Value = "covid"
$ws.Cells.Item(8,3).Value = 0.5154443338717717
$ws.Cells.Item(8,6).ClearContents()
$ws.Cells.Item(9,2).Value = "covid pcr"
$ws.Cells.Item(9,3).Value = 0.585422678668393
$ws.Cells.Item(9,6).Value = 0.6389300840389701
$ws.Cells.Item(10,2).Value = "pcr exame covid"
$ws.Cells.Item(10,3).Value = 0.6213586636848555
$ws.Cells.Item(10,6).Value = 0.6489958741777089
$ws.Cells.Item(11,2).Value = "pcr"
$ws.Cells.Item(11,3).Value = 0.6122924582820127
$ws.Cells.Item(11,6).Value = 0.6850871894141481
$ws.Cells.Item(12,2).Value = "decreto lockdown"
$ws.Cells.Item(12,3).Value = 0.638803966281798
$ws.Cells.Item(12,5).ClearContents()
$ws.Cells.Item(13,2).Value = "estou com covid"
$ws.Cells.Item(13,3).Value = 0.6346744245229214
$ws.Cells.Item(13,4).Value = 0.6981139703447465
$ws.Cells.Item(13,5).Value = 0.5557749936627187
$ws.Cells.Item(13,6).ClearContents()
$ws.Cells.Item(14,2).Value = "peguei covid"
$ws.Cells.Item(14,3).Value = 0.5891242098258881
$ws.Cells.Item(14,4).ClearContents()
$ws.Cells.Item(14,5).ClearContents()
$ws.Cells.Item(14,6).Value = 0.5764232940728571
$ws.Cells.Item(15,2).Value = "febre"
$ws.Cells.Item(15,4).Value = -0.5331219916139754
$ws.Cells.Item(15,5).Value = -0.6686306366285906
$ws.Cells.Item(15,6).ClearContents()
$ws.Cells.Item(16,2).Value = "brasil coronavírus"
$ws.Cells.Item(16,4).Value = -0.5325766412243701
$ws.Cells.Item(16,5).Value = -0.6778293702356335
$ws.Cells.Item(16,6).Value = -0.6145569668541886
$ws.Cells.Item(17,2).Value = "coronavírus no brasil"
$ws.Cells.Item(17,4).Value = -0.5475383839781843
$ws.Cells.Item(17,5).Value = -0.6975044606138195
$ws.Cells.Item(17,6).Value = -0.6385059533175508
$ws.Cells.Item(18,2).Value = "corona grupo de risco"
$ws.Cells.Item(18,4).Value = -0.5844662603638663
$ws.Cells.Item(18,5).Value = -0.6818408783946586
$ws.Cells.Item(18,6).Value = -0.6837977640347525
$ws.Cells.Item(19,2).Value = "ministerio da saude"
$ws.Cells.Item(19,4).Value = -0.5730336781793793
$ws.Cells.Item(19,6).Value = -0.6572906881215554
$ws.Cells.Item(20,2).Value = "taxa de ocupação de leitos"
$ws.Cells.Item(20,4).Value = 0.6134354751192254
$ws.Cells.Item(20,5).ClearContents()
$ws.Cells.Item(20,6).Value = 0.5025116353919563
$ws.Cells.Item(21,2).Value = "teste rápido covid"
$ws.Cells.Item(21,4).Value = 0.6399211101609452
$ws.Cells.Item(21,5).Value = 0.5048475763776782
$ws.Cells.Item(22,2).Value = "teste igg"
$ws.Cells.Item(22,4).Value = 0.6909925168547132
$ws.Cells.Item(22,5).Value = 0.5201312155230962
$ws.Cells.Item(22,6).ClearContents()
$ws.Cells.Item(23,2).Value = "exame igg"
$ws.Cells.Item(23,4).Value = 0.6049923399445443
$ws.Cells.Item(24,2).Value = "mortes corona"
$ws.Cells.Item(24,4).Value = -0.5492235398299186
$ws.Cells.Item(24,5).Value = -0.6897909803467545
$ws.Cells.Item(24,6).Value = -0.6471259579868259
$ws.Cells.Item(25,2).Value = "coronavírus oms"
$ws.Cells.Item(25,4).Value = -0.5076437718615632
$ws.Cells.Item(25,5).Value = -0.6101702646448129
$ws.Cells.Item(25,6).Value = -0.5619736088228765
$ws.Cells.Item(26,2).Value = "oms corona"
$ws.Cells.Item(26,4).Value = -0.5879917759284979
$ws.Cells.Item(26,5).ClearContents()
$ws.Cells.Item(27,2).Value = "oms coronavírus"
$ws.Cells.Item(27,4).Value = -0.5076437718615632
$ws.Cells.Item(27,5).Value = -0.6101702646448129
$ws.Cells.Item(27,6).Value = -0.5619736088228765
$ws.Cells.Item(28,2).Value = "álcool gel"
$ws.Cells.Item(28,4).Value = -0.5614336830723569
$ws.Cells.Item(28,5).Value = -0.6808175897129142
$ws.Cells.Item(28,6).Value = -0.672871239563938
$ws.Cells.Item(29,2).Value = "álcool 70"
$ws.Cells.Item(29,4).Value = -0.5885393116760118
$ws.Cells.Item(29,5).ClearContents()
$ws.Cells.Item(30,2).Value = "respirador mecanico"
$ws.Cells.Item(30,4).Value = -0.557995069952795
$ws.Cells.Item(30,5).Value = -0.6333132793984371
$ws.Cells.Item(30,6).Value = -0.5472225015598507
$ws.Cells.Item(31,2).Value = "medicamento corona"
$ws.Cells.Item(31,4).Value = -0.6267466041924911
$ws.Cells.Item(31,5).ClearContents()
$ws.Cells.Item(32,2).Value = "teste covid"
$ws.Cells.Item(32,5).Value = 0.6153846704505471
$ws.Cells.Item(33,2).Value = "tosse"
$ws.Cells.Item(33,5).Value = -0.5589999006602973
$ws.Cells.Item(33,6).Value = -0.6855960655591937
$ws.Cells.Item(34,2).Value = "sintomas coronavirus"
$ws.Cells.Item(34,5).Value = -0.6522066767297013
$ws.Cells.Item(34,6).Value = -0.6950522024227365
$ws.Cells.Item(35,2).Value = "tomar ivermectina"
$ws.Cells.Item(35,5).Value = 0.6744555319420698
$ws.Cells.Item(35,6).Value = 0.5734724143898801
$ws.Cells.Item(36,5).Value = -0.5671498816763793
$ws.Cells.Item(36,6).Value = -0.5766634751667631
$ws.Cells.Item(37,2).Value = "máscara n95"
$ws.Cells.Item(37,5).Value = -0.5480085643904264
$ws.Cells.Item(37,6).Value = -0.6374713944119883
$ws.Cells.Item(38,2).Value = "máscara descartável"
$ws.Cells.Item(38,5).Value = -0.5274680874427935
$ws.Cells.Item(38,6).Value = -0.5819727646862813
$ws.Cells.Item(39,2).Value = "igg"
$ws.Cells.Item(39,5).Value = 0.5929138644916055
$ws.Cells.Item(39,6).ClearContents()
$ws.Cells.Item(40,2).Value = "igm"
$ws.Cells.Item(40,5).Value = 0.5989424214441603
$ws.Cells.Item(40,6).ClearContents()
$ws.Cells.Item(41,2).Value = "igg igm"
$ws.Cells.Item(41,5).Value = 0.5747471120066769
$ws.Cells.Item(42,2).Value = "covid igg"
$ws.Cells.Item(42,5).Value = 0.5746413627680371
$ws.Cells.Item(43,2).Value = "covid igm igg"
$ws.Cells.Item(43,5).Value = 0.5783184108169995
$ws.Cells.Item(43,6).ClearContents()
$ws.Cells.Item(44,2).Value = "coronavírus quarentena"
$ws.Cells.Item(44,5).Value = -0.5752735762357104
$ws.Cells.Item(44,6).Value = -0.5633235210016591
$ws.Cells.Item(45,2).Value = "quarentena brasil"
$ws.Cells.Item(45,5).Value = -0.5036670157482387
$ws.Cells.Item(46,2).Value = "covid 19 oms"
$ws.Cells.Item(46,5).Value = -0.5051023846640987
$ws.Cells.Item(47,2).Value = "oms brasil coronavírus"
$ws.Cells.Item(47,5).Value = -0.5148097968856272
$ws.Cells.Item(47,6).Value = -0.5419616172372991
$ws.Cells.Item(48,2).Value = "coronavírus imunidade"
$ws.Cells.Item(48,5).Value = -0.566463536284767
$ws.Cells.Item(48,6).Value = -0.5553794830693396
$ws.Cells.Item(49,2).Value = "respirador"
$ws.Cells.Item(49,5).Value = -0.5964249241875299
$ws.Cells.Item(49,6).Value = -0.5439442563443259
$ws.Cells.Item(50,2).Value = "medicamento coronavírus"
$ws.Cells.Item(50,5).Value = -0.5203046143040917
$ws.Cells.Item(50,6).Value = -0.5397565412264164
$ws.Cells.Item(51,2).Value = "dor de garganta"
$ws.Cells.Item(51,5).ClearContents()
$ws.Cells.Item(51,6).Value = -0.553773984610456
$ws.Cells.Item(52,2).Value = "exame cotonete"
$ws.Cells.Item(52,5).ClearContents()
$ws.Cells.Item(52,6).Value = 0.6704883530278225
$ws.Cells.Item(53,2).Value = "teste coronavírus"
$ws.Cells.Item(53,5).ClearContents()
$ws.Cells.Item(53,6).Value = -0.5063532414140993
$ws.Cells.Item(54,2).Value = "ivermectina covid como tomar"
$ws.Cells.Item(54,5).ClearContents()
$ws.Cells.Item(54,6).Value = 0.5998381010958729
$ws.Cells.Item(55,2).Value = "coronavírus grupo de risco"
$ws.Cells.Item(55,6).Value = -0.5193770984603348
$ws.Cells.Item(56,2).Value = "vacina corona"
$ws.Cells.Item(56,6).Value = -0.5730748275408113
$ws.Cells.Item(57,2).Value = "exame covid"
$ws.Cells.Item(57,6).Value = 0.548867454461316
$ws.Cells.Item(58,2).Value = "reagente igg"
$ws.Cells.Item(58,6).Value = 0.6884460304254997
$ws.Cells.Item(59,2).Value = "respirador hospitalar"
$ws.Cells.Item(59,6).Value = -0.5708894331659432
